$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell holds plain text in the source data (coin prices/volumes are
# formatted strings, not numbers). Writing with a leading apostrophe forces
# Excel to store the literal text instead of auto-converting look-alikes
# such as "1.00" or "8.30" into numbers; resetting the style back to Normal
# afterwards avoids leaving a stray "Text" number-format style on the cell.

# Row 2
$ws.Range("D2").Value = "'94.685.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.40%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.523.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.06%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'239.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.35%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'631.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.73%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +5.23%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.397"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.96%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  -0.09%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +7.82%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'3.517.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.87%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'43.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.13%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  +4.86%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'6.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.78%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'4.186.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.14%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'94.543.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.49%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "'  +3.67%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'8.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.93%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'3.521.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.93%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'12.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +15.26%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'18.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.68%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.501"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +9.83%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'518.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +5.59%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  +0.56%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'6.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +8.09%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  +3.09%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'92.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.32%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'12.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.32%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'2.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +11.19%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'11.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.29%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'0.142"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.81%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  +0.01%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'0.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.13%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.22%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("B35").Value = "'EthereumClassic"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'30.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.64%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("B36").Value = "'PolygonEcosystemToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.562"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.20%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'586.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +10.51%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  +6.52%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'7.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.84%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  +0.01%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "'ARBITRUM"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.928"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.84%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'Kaspa"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.150"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.93%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'23.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.16%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.0423"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.30%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'1.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.38%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'5.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.22%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  -0.15%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'  +1.98%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'54.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.28%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'8.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.19%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  +1.08%  "
$ws.Range("E51").Style = "Normal"

